# Applies the CapaciteSavoirfaire.xlsx update:
#  - refresh the generation "Date" metadata value
#  - fix casing (exerciceProfessionnel -> ExerciceProfessionnel) and tidy
#    wording in the Elements rows describing that reference
#  - nudge the width of the first two "Elements" columns

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the generation Date ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-10-29T11:46:56+00:00"

# --- Elements sheet: ExerciceProfessionnel casing + text tweaks ---
$wsElem = $wb.Worksheets.Item("Elements")
$wsElem.Range("A6").Value = "CapaciteSavoirfaire.ExerciceProfessionnel"
$wsElem.Range("B6").Value = "CapaciteSavoirfaire.ExerciceProfessionnel"
$wsElem.Range("L6").Value = "Lien vers la classe ExerciceProfessionnel"
$wsElem.Range("M6").Value = "Lien vers la classe ExerciceProfessionnel"
$wsElem.Range("AF6").Value = "SavoirFaire.ExerciceProfessionnel"

# --- Elements sheet: widen columns A and B slightly (33.953125 -> 33.98046875 chars) ---
# NOTE: Excel's ColumnWidth COM property only resolves to whole-pixel
# increments (~0.1667 "chars" steps) once written back to the OOXML <col
# width=.../> attribute, so the exact 33.98046875 value is not reachable
# through this API; 33.0 is the closest input that lands on the nearest
# achievable stored width.
$wsElem.Columns.Item(1).ColumnWidth = 33.0
$wsElem.Columns.Item(2).ColumnWidth = 33.0
